$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''293.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''-6.91%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''40.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''-1.21%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.033'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-2.03%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.07326'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''-3.64%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''1.537'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-8.62%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.9293'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''-0.09%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = '''0.1168'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-2.85%'
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.1740'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''-4.44%'
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.04347'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''5.11%'
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.08735'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-2.98%'
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.1054'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''-0.12%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = '''BitForexToken'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''0.001272'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-0.65%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = '''CoinExToken'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '''https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''0.03942'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-0.69%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = '''TigerCash'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''0.005937'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''0.59%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = '''LEO'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''3.337'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''0.05%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = '''GateToken'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '''4.283'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-0.79%'
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = '''BitpandaEcosystemToken'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '''https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '''0.3289'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''-1.76%'
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = '''MCDex'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '''https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '''7.976'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''4.68%'
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = '''ProBitToken'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = '''https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = '''0.1400'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''3.80%'
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = '''ZBToken'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '''https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = '''0.2742'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''-3.17%'
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''-1.35%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.003787'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '''-5.14%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0003725'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''22.56%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = '''0.02303'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''-5.10%'
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.05072'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''-1.88%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.006225'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''88.66%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.007868'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''1.70%'
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''-1.18%'
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.007377'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''-2.70%'
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''-3.43%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.3190'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-6.01%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00006279'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''-4.72%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.00000000751'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''0.10%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.03291'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''-87.82%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.00002102'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''0.10%'
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.0002002'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''0.10%'
$ws.Range("E50").Style = "Normal"
